$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values before overwriting anything
$a1 = $ws.Range("A1").Value()
$b1 = $ws.Range("B1").Value()

$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()

$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()

$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()

# Swap header labels: column A becomes "Dept No", column B becomes "Dept Name"
$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

# Swap data columns A/B for each department row (Dept No moves to A, Dept Name moves to B)
$ws.Range("A2").Value = $b2
$ws.Range("B2").Value = $a2

$ws.Range("A3").Value = $b3
$ws.Range("B3").Value = $a3

$ws.Range("A4").Value = $b4
$ws.Range("B4").Value = $a4

# Update Total_Compensation values in column C
$ws.Range("C2").Value = 4299600
$ws.Range("C4").Value = 4657650
